# DialogInheritance.xlsx - add the "popAbwesenheitAnlegen" dialog row
# (msz - restructuring control processing -> container)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New data row: dialog class name in column A, window title in column D
# (mirrors the existing dlgProfil / "Mein Profil | TT-Planer" row above it)
$ws.Range("A6").Value = "popAbwesenheitAnlegen"
$ws.Range("D6").Value = "Abwesenheiten | TT-Planer"

# Column D needs to widen (best-fit) to show the new, longer title text
$ws.Columns.Item(4).ColumnWidth = 22.39

# Leave the selection on the freshly-entered cell, like the author did
$ws.Range("A6").Select()

# Page setup was touched as well (paper size / orientation)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
